# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity, Temperature, and mmWave sheets,
# matching new telemetry captured on 2026-01-28 (17:19:40 - 17:20:18 local time).

$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 153-162 (motion sensor, Bathroom) ---
$ws = $wb.Worksheets.Item("PIR")
$pirData = @(
    ,@(153, "2026-01-28", "17:19:40", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(154, "2026-01-28", "17:19:41", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(155, "2026-01-28", "17:19:43", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(156, "2026-01-28", "17:19:48", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(157, "2026-01-28", "17:19:54", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(158, "2026-01-28", "17:19:58", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(159, "2026-01-28", "17:20:03", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(160, "2026-01-28", "17:20:08", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(161, "2026-01-28", "17:20:14", "17:00", "Bathroom", "No Motion", "Inactive")
    ,@(162, "2026-01-28", "17:20:18", "17:00", "Bathroom", "No Motion", "Inactive")
)
foreach ($row in $pirData) {
    $r = $row[0]
    $aCell = $ws.Range("A$r")
    $aCell.NumberFormat = "@"
    $aCell.Value = $row[1]
    $aCell.ClearFormats()
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}

# --- Humidity sheet: rows 153-162 (Bathroom, percentage values need forced text) ---
$ws = $wb.Worksheets.Item("Humidity")
$humidityData = @(
    ,@(153, "2026-01-28", "17:19:40", "17:00", "Bathroom", "87.4%", "Active")
    ,@(154, "2026-01-28", "17:19:42", "17:00", "Bathroom", "86.6%", "Active")
    ,@(155, "2026-01-28", "17:19:44", "17:00", "Bathroom", "87.5%", "Active")
    ,@(156, "2026-01-28", "17:19:49", "17:00", "Bathroom", "87.5%", "Active")
    ,@(157, "2026-01-28", "17:19:52", "17:00", "Bathroom", "86.6%", "Active")
    ,@(158, "2026-01-28", "17:19:56", "17:00", "Bathroom", "87.5%", "Active")
    ,@(159, "2026-01-28", "17:20:05", "17:00", "Bathroom", "86.6%", "Active")
    ,@(160, "2026-01-28", "17:20:09", "17:00", "Bathroom", "87.5%", "Active")
    ,@(161, "2026-01-28", "17:20:13", "17:00", "Bathroom", "86.6%", "Active")
    ,@(162, "2026-01-28", "17:20:17", "17:00", "Bathroom", "86.1%", "Active")
)
foreach ($row in $humidityData) {
    $r = $row[0]
    $aCell = $ws.Range("A$r")
    $aCell.NumberFormat = "@"
    $aCell.Value = $row[1]
    $aCell.ClearFormats()
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $eCell = $ws.Range("E$r")
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[5]
    $eCell.ClearFormats()
    $ws.Range("F$r").Value = $row[6]
}

# --- Temperature sheet: rows 153-162 (Bathroom) ---
$ws = $wb.Worksheets.Item("Temperature")
$temperatureData = @(
    ,@(153, "2026-01-28", "17:19:41", "17:00", "Bathroom", "22.8C", "Active")
    ,@(154, "2026-01-28", "17:19:43", "17:00", "Bathroom", "22.8C", "Active")
    ,@(155, "2026-01-28", "17:19:45", "17:00", "Bathroom", "22.8C", "Active")
    ,@(156, "2026-01-28", "17:19:49", "17:00", "Bathroom", "22.8C", "Active")
    ,@(157, "2026-01-28", "17:19:53", "17:00", "Bathroom", "22.8C", "Active")
    ,@(158, "2026-01-28", "17:19:57", "17:00", "Bathroom", "22.8C", "Active")
    ,@(159, "2026-01-28", "17:20:05", "17:00", "Bathroom", "22.8C", "Active")
    ,@(160, "2026-01-28", "17:20:09", "17:00", "Bathroom", "22.8C", "Active")
    ,@(161, "2026-01-28", "17:20:13", "17:00", "Bathroom", "22.8C", "Active")
    ,@(162, "2026-01-28", "17:20:17", "17:00", "Bathroom", "22.8C", "Active")
)
foreach ($row in $temperatureData) {
    $r = $row[0]
    $aCell = $ws.Range("A$r")
    $aCell.NumberFormat = "@"
    $aCell.Value = $row[1]
    $aCell.ClearFormats()
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}

# --- mmWave sheet: rows 31-43 (Living Room presence sensor) ---
$ws = $wb.Worksheets.Item("mmWave")
$mmwaveData = @(
    ,@(31, "2026-01-28", "17:19:41", "17:00", "Living Room", "NO_PRESENCE", "Inactive")
    ,@(32, "2026-01-28", "17:19:42", "17:00", "Living Room", "NO_PRESENCE", "Inactive")
    ,@(33, "2026-01-28", "17:19:44", "17:00", "Living Room", "NO_PRESENCE", "Inactive")
    ,@(34, "2026-01-28", "17:19:46", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(35, "2026-01-28", "17:19:50", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(36, "2026-01-28", "17:19:53", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(37, "2026-01-28", "17:19:55", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(38, "2026-01-28", "17:19:59", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(39, "2026-01-28", "17:20:02", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(40, "2026-01-28", "17:20:04", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(41, "2026-01-28", "17:20:07", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(42, "2026-01-28", "17:20:10", "17:00", "Living Room", "PRESENCE", "Active")
    ,@(43, "2026-01-28", "17:20:14", "17:00", "Living Room", "PRESENCE", "Active")
)
foreach ($row in $mmwaveData) {
    $r = $row[0]
    $aCell = $ws.Range("A$r")
    $aCell.NumberFormat = "@"
    $aCell.Value = $row[1]
    $aCell.ClearFormats()
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}

